$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CHL-A")

# Update sample IDs in column A (refinements to chl a data work up)
$ws.Range("A36").Value = "07262016.WKA.SU-30.3.1.UNK"
$ws.Range("A37").Value = "07262016.WKA.SU-30.3.1.UNK-A"
$ws.Range("A38").Value = "07262016.WKA.SU-01.15.UNK"
$ws.Range("A39").Value = "07262016.WKA.SU-01.15.UNK-A"

# Reflect the user's final scroll position / selection in the saved view
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 22
$win.ScrollColumn = 1
$ws.Range("A39").Select() | Out-Null
